$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 6779.3335
$ws.Cells.Item(17, 10).Value = 6128.125
$ws.Cells.Item(17, 12).Value = 18384.375
$ws.Cells.Item(17, 14).Value = -18720.375
$ws.Cells.Item(18, 8).Value = 1001
$ws.Cells.Item(18, 9).Value = 1001
$ws.Cells.Item(18, 11).Value = 1001
$ws.Cells.Item(18, 13).Value = -717
$ws.Cells.Item(19, 8).Value = 4054.348
$ws.Cells.Item(19, 9).Value = 3501.7693
$ws.Cells.Item(19, 11).Value = 3501.7693
$ws.Cells.Item(19, 13).Value = -3326.7693
$ws.Cells.Item(33, 8).Value = 292.25
$ws.Cells.Item(33, 9).Value = 281.86365
$ws.Cells.Item(33, 11).Value = 281.86365
$ws.Cells.Item(33, 13).Value = -52.86365000000001
$ws.Cells.Item(69, 8).Value = 19291.889
$ws.Cells.Item(69, 9).Value = 18197
$ws.Cells.Item(69, 10).Value = 19713
$ws.Cells.Item(69, 11).Value = 54591
$ws.Cells.Item(69, 12).Value = 59139
$ws.Cells.Item(69, 13).Value = -53717
$ws.Cells.Item(69, 14).Value = -60887
$ws.Cells.Item(70, 8).Value = 1664191.5
$ws.Cells.Item(70, 9).Value = 2113.8
$ws.Cells.Item(70, 10).Value = 2587568
$ws.Cells.Item(70, 11).Value = 6341.400000000001
$ws.Cells.Item(70, 12).Value = 7762704
$ws.Cells.Item(70, 13).Value = -6071.400000000001
$ws.Cells.Item(70, 14).Value = -7763244
$ws.Cells.Item(72, 8).Value = 19291.889
$ws.Cells.Item(72, 9).Value = 18197
$ws.Cells.Item(72, 10).Value = 19713
$ws.Cells.Item(72, 11).Value = 163773
$ws.Cells.Item(72, 12).Value = 177417
$ws.Cells.Item(72, 13).Value = -159405
$ws.Cells.Item(72, 14).Value = -186153
$ws.Cells.Item(73, 8).Value = 1664191.5
$ws.Cells.Item(73, 9).Value = 2113.8
$ws.Cells.Item(73, 10).Value = 2587568
$ws.Cells.Item(73, 11).Value = 6341.400000000001
$ws.Cells.Item(73, 12).Value = 7762704
$ws.Cells.Item(73, 13).Value = -5405.400000000001
$ws.Cells.Item(73, 14).Value = -7764576
$ws.Cells.Item(80, 8).Value = 1333.875
$ws.Cells.Item(80, 9).Value = 302.5
$ws.Cells.Item(80, 10).Value = 1952.7
$ws.Cells.Item(80, 11).Value = 907.5
$ws.Cells.Item(80, 12).Value = 5858.1
$ws.Cells.Item(80, 13).Value = 90.5
$ws.Cells.Item(80, 14).Value = -7854.1
$ws.Cells.Item(83, 8).Value = 1333.875
$ws.Cells.Item(83, 9).Value = 302.5
$ws.Cells.Item(83, 10).Value = 1952.7
$ws.Cells.Item(83, 11).Value = 2722.5
$ws.Cells.Item(83, 12).Value = 17574.3
$ws.Cells.Item(83, 13).Value = 2269.5
$ws.Cells.Item(83, 14).Value = -27558.3
$ws.Cells.Item(100, 8).Value = 9970.75
$ws.Cells.Item(100, 9).Value = 3288.6667
$ws.Cells.Item(100, 10).Value = 12198.111
$ws.Cells.Item(100, 11).Value = 3288.6667
$ws.Cells.Item(100, 12).Value = 12198.111
$ws.Cells.Item(100, 13).Value = -2747.6667
$ws.Cells.Item(100, 14).Value = -13280.111
$ws.Cells.Item(106, 8).Value = 3562.5
$ws.Cells.Item(106, 10).Value = 3562.5
$ws.Cells.Item(106, 12).Value = 3562.5
$ws.Cells.Item(106, 14).Value = -4824.5
$ws.Cells.Item(113, 8).Value = 168667.17
$ws.Cells.Item(113, 9).Value = 2499.75
$ws.Cells.Item(113, 11).Value = 2499.75
$ws.Cells.Item(113, 13).Value = 754.25
$ws.Cells.Item(121, 8).Value = 1366.3334
$ws.Cells.Item(121, 10).Value = 1366.3334
$ws.Cells.Item(121, 12).Value = 4099.0002
$ws.Cells.Item(121, 14).Value = -7593.0002
$ws.Cells.Item(125, 8).Value = 6499
$ws.Cells.Item(125, 9).Value = 6499
$ws.Cells.Item(125, 11).Value = 58491
$ws.Cells.Item(125, 13).Value = -56031
$ws.Cells.Item(129, 8).Value = 2299.923
$ws.Cells.Item(129, 9).Value = 1988
$ws.Cells.Item(129, 10).Value = 2393.5
$ws.Cells.Item(129, 11).Value = 5964
$ws.Cells.Item(129, 12).Value = 7180.5
$ws.Cells.Item(129, 13).Value = -964
$ws.Cells.Item(129, 14).Value = -17180.5
$ws.Cells.Item(131, 8).Value = 6671.1
$ws.Cells.Item(131, 9).Value = 3186.0833
$ws.Cells.Item(131, 10).Value = 8994.444
$ws.Cells.Item(131, 11).Value = 9558.249899999999
$ws.Cells.Item(131, 12).Value = 26983.332
$ws.Cells.Item(131, 13).Value = -4518.249899999999
$ws.Cells.Item(131, 14).Value = -37063.33199999999
$ws.Cells.Item(132, 8).Value = 1964.0385
$ws.Cells.Item(132, 9).Value = 1312.8572
$ws.Cells.Item(132, 11).Value = 3938.5716
$ws.Cells.Item(132, 13).Value = -1408.5716
$ws.Cells.Item(135, 8).Value = 2171.8096
$ws.Cells.Item(135, 9).Value = 1440.5333
$ws.Cells.Item(135, 10).Value = 4000
$ws.Cells.Item(135, 11).Value = 12964.7997
$ws.Cells.Item(135, 12).Value = 36000
$ws.Cells.Item(135, 13).Value = -10429.7997
$ws.Cells.Item(135, 14).Value = -41070
$ws.Cells.Item(137, 8).Value = 7838.409
$ws.Cells.Item(137, 9).Value = 3097.125
$ws.Cells.Item(137, 10).Value = 10547.714
$ws.Cells.Item(137, 11).Value = 9291.375
$ws.Cells.Item(137, 12).Value = 31643.142
$ws.Cells.Item(137, 13).Value = -6741.375
$ws.Cells.Item(137, 14).Value = -36743.142
$ws.Cells.Item(138, 8).Value = 4515.73
$ws.Cells.Item(138, 10).Value = 5000.6313
$ws.Cells.Item(138, 12).Value = 15001.8939
$ws.Cells.Item(138, 14).Value = -25281.8939
$ws.Cells.Item(141, 8).Value = 1539.6552
$ws.Cells.Item(141, 9).Value = 1487.3214
$ws.Cells.Item(141, 10).Value = 3005
$ws.Cells.Item(141, 11).Value = 4461.9642
$ws.Cells.Item(141, 12).Value = 9015
$ws.Cells.Item(141, 13).Value = 718.0357999999997
$ws.Cells.Item(141, 14).Value = -19375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 25135.111
$ws.Cells.Item(32, 9).Value = 23805.555
$ws.Cells.Item(32, 10).Value = 34442
$ws.Cells.Item(32, 11).Value = 23805.555
$ws.Cells.Item(32, 12).Value = 34442
$ws.Cells.Item(32, 13).Value = -23518.555
$ws.Cells.Item(32, 14).Value = -35016
$ws.Cells.Item(38, 8).Value = 3000
$ws.Cells.Item(38, 9).Value = 3000
$ws.Cells.Item(38, 11).Value = 3000
$ws.Cells.Item(38, 13).Value = -2533
$ws.Cells.Item(45, 8).Value = 2199.8
$ws.Cells.Item(45, 10).Value = 3999.5
$ws.Cells.Item(45, 12).Value = 3999.5
$ws.Cells.Item(45, 14).Value = -4753.5
$ws.Cells.Item(61, 8).Value = 11007.241
$ws.Cells.Item(61, 9).Value = 7351.222
$ws.Cells.Item(61, 10).Value = 16989.818
$ws.Cells.Item(61, 11).Value = 7351.222
$ws.Cells.Item(61, 12).Value = 16989.818
$ws.Cells.Item(61, 13).Value = -7139.222
$ws.Cells.Item(61, 14).Value = -17413.818
$ws.Cells.Item(74, 8).Value = 716568.6
$ws.Cells.Item(74, 9).Value = 835330.0600000001
$ws.Cells.Item(74, 10).Value = 4000
$ws.Cells.Item(74, 11).Value = 835330.0600000001
$ws.Cells.Item(74, 12).Value = 4000
$ws.Cells.Item(74, 13).Value = -834456.0600000001
$ws.Cells.Item(74, 14).Value = -5748
$ws.Cells.Item(77, 8).Value = 716568.6
$ws.Cells.Item(77, 9).Value = 835330.0600000001
$ws.Cells.Item(77, 10).Value = 4000
$ws.Cells.Item(77, 11).Value = 4176650.3
$ws.Cells.Item(77, 12).Value = 20000
$ws.Cells.Item(77, 13).Value = -4172282.3
$ws.Cells.Item(77, 14).Value = -28736
$ws.Cells.Item(122, 8).Value = 4635.3335
$ws.Cells.Item(122, 9).Value = 4053.8667
$ws.Cells.Item(122, 11).Value = 12161.6001
$ws.Cells.Item(122, 13).Value = -9711.6001
$ws.Cells.Item(132, 8).Value = 6801.698
$ws.Cells.Item(132, 9).Value = 2802.1914
$ws.Cells.Item(132, 10).Value = 18550.25
$ws.Cells.Item(132, 11).Value = 8406.574200000001
$ws.Cells.Item(132, 12).Value = 55650.75
$ws.Cells.Item(132, 13).Value = -5876.574200000001
$ws.Cells.Item(132, 14).Value = -60710.75
$ws.Cells.Item(136, 8).Value = 11007.241
$ws.Cells.Item(136, 9).Value = 7351.222
$ws.Cells.Item(136, 10).Value = 16989.818
$ws.Cells.Item(136, 11).Value = 22053.666
$ws.Cells.Item(136, 12).Value = 50969.454
$ws.Cells.Item(136, 13).Value = -19503.666
$ws.Cells.Item(136, 14).Value = -56069.454

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 41502.5
$ws.Cells.Item(20, 9).Value = 58396.668
$ws.Cells.Item(20, 10).Value = 3490.625
$ws.Cells.Item(20, 11).Value = 58396.668
$ws.Cells.Item(20, 12).Value = 3490.625
$ws.Cells.Item(20, 13).Value = -58149.668
$ws.Cells.Item(20, 14).Value = -3984.625
$ws.Cells.Item(22, 8).Value = 1965
$ws.Cells.Item(22, 9).Value = 2111.375
$ws.Cells.Item(22, 10).Value = 794
$ws.Cells.Item(22, 11).Value = 2111.375
$ws.Cells.Item(22, 12).Value = 794
$ws.Cells.Item(22, 13).Value = -1938.375
$ws.Cells.Item(22, 14).Value = -1140
$ws.Cells.Item(23, 8).Value = 1014
$ws.Cells.Item(23, 10).Value = 1014
$ws.Cells.Item(23, 12).Value = 1014
$ws.Cells.Item(23, 14).Value = -1580
$ws.Cells.Item(86, 8).Value = 88712.25999999999
$ws.Cells.Item(86, 9).Value = 1879.6842
$ws.Cells.Item(86, 10).Value = 501167
$ws.Cells.Item(86, 11).Value = 1879.6842
$ws.Cells.Item(86, 12).Value = 501167
$ws.Cells.Item(86, 13).Value = -756.6841999999999
$ws.Cells.Item(86, 14).Value = -503413
$ws.Cells.Item(89, 8).Value = 88712.25999999999
$ws.Cells.Item(89, 9).Value = 1879.6842
$ws.Cells.Item(89, 10).Value = 501167
$ws.Cells.Item(89, 11).Value = 9398.421
$ws.Cells.Item(89, 12).Value = 2505835
$ws.Cells.Item(89, 13).Value = -3782.421
$ws.Cells.Item(89, 14).Value = -2517067
$ws.Cells.Item(99, 8).Value = 2507.1667
$ws.Cells.Item(99, 9).Value = 2213.7368
$ws.Cells.Item(99, 10).Value = 3622.2
$ws.Cells.Item(99, 11).Value = 2213.7368
$ws.Cells.Item(99, 12).Value = 3622.2
$ws.Cells.Item(99, 13).Value = -715.7368000000001
$ws.Cells.Item(99, 14).Value = -6618.2
$ws.Cells.Item(105, 8).Value = 55572790
$ws.Cells.Item(105, 9).Value = 90934770
$ws.Cells.Item(105, 10).Value = 3968.7144
$ws.Cells.Item(105, 11).Value = 90934770
$ws.Cells.Item(105, 12).Value = 3968.7144
$ws.Cells.Item(105, 13).Value = -90933023
$ws.Cells.Item(105, 14).Value = -7462.7144
$ws.Cells.Item(107, 8).Value = 1854.48
$ws.Cells.Item(107, 9).Value = 1758.375
$ws.Cells.Item(107, 11).Value = 1758.375
$ws.Cells.Item(107, 13).Value = 161.625
$ws.Cells.Item(134, 8).Value = 4309.78
$ws.Cells.Item(134, 9).Value = 3536.3333
$ws.Cells.Item(134, 10).Value = 6795.857
$ws.Cells.Item(134, 11).Value = 10608.9999
$ws.Cells.Item(134, 12).Value = 20387.571
$ws.Cells.Item(134, 13).Value = -8073.999899999999
$ws.Cells.Item(134, 14).Value = -25457.571
$ws.Cells.Item(141, 8).Value = 91666.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(5, 8).Value = 304.7143
$ws.Cells.Item(5, 9).Value = 311.4
$ws.Cells.Item(5, 11).Value = 311.4
$ws.Cells.Item(5, 13).Value = -199.4
$ws.Cells.Item(6, 8).Value = 6833.3335
$ws.Cells.Item(6, 9).Value = 7500.5
$ws.Cells.Item(6, 11).Value = 7500.5
$ws.Cells.Item(6, 13).Value = -7387.5
$ws.Cells.Item(7, 8).Value = 368.42307
$ws.Cells.Item(7, 10).Value = 129.75
$ws.Cells.Item(7, 12).Value = 129.75
$ws.Cells.Item(7, 14).Value = -355.75
$ws.Cells.Item(10, 8).Value = 2977.5
$ws.Cells.Item(10, 10).Value = 5504
$ws.Cells.Item(10, 12).Value = 5504
$ws.Cells.Item(10, 14).Value = -5782
$ws.Cells.Item(29, 8).Value = 10000
$ws.Cells.Item(29, 10).Value = 10000
$ws.Cells.Item(29, 12).Value = 10000
$ws.Cells.Item(29, 14).Value = -10586
$ws.Cells.Item(31, 8).Value = 17863542
$ws.Cells.Item(31, 9).Value = 62503116
$ws.Cells.Item(31, 10).Value = 7713.225
$ws.Cells.Item(31, 11).Value = 62503116
$ws.Cells.Item(31, 12).Value = 7713.225
$ws.Cells.Item(31, 13).Value = -62502821
$ws.Cells.Item(31, 14).Value = -8303.225
$ws.Cells.Item(32, 8).Value = 7500
$ws.Cells.Item(32, 9).Value = 7500
$ws.Cells.Item(32, 11).Value = 7500
$ws.Cells.Item(32, 13).Value = -7184
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 14).ClearContents()
$ws.Cells.Item(34, 8).Value = 17863542
$ws.Cells.Item(34, 9).Value = 62503116
$ws.Cells.Item(34, 10).Value = 7713.225
$ws.Cells.Item(34, 11).Value = 62503116
$ws.Cells.Item(34, 12).Value = 7713.225
$ws.Cells.Item(34, 13).Value = -62502914
$ws.Cells.Item(34, 14).Value = -8117.225
$ws.Cells.Item(38, 8).Value = 8082.375
$ws.Cells.Item(38, 9).Value = 12654.25
$ws.Cells.Item(38, 10).Value = 3510.5
$ws.Cells.Item(38, 11).Value = 12654.25
$ws.Cells.Item(38, 12).Value = 3510.5
$ws.Cells.Item(38, 13).Value = -12277.25
$ws.Cells.Item(38, 14).Value = -4264.5
$ws.Cells.Item(46, 8).Value = 8082.375
$ws.Cells.Item(46, 9).Value = 12654.25
$ws.Cells.Item(46, 10).Value = 3510.5
$ws.Cells.Item(46, 11).Value = 12654.25
$ws.Cells.Item(46, 12).Value = 3510.5
$ws.Cells.Item(46, 13).Value = -12443.25
$ws.Cells.Item(46, 14).Value = -3932.5
$ws.Cells.Item(58, 8).Value = 4549
$ws.Cells.Item(58, 9).Value = 3375.5454
$ws.Cells.Item(58, 10).Value = 6162.5
$ws.Cells.Item(58, 11).Value = 3375.5454
$ws.Cells.Item(58, 12).Value = 6162.5
$ws.Cells.Item(58, 13).Value = -3172.5454
$ws.Cells.Item(58, 14).Value = -6568.5
$ws.Cells.Item(86, 8).Value = 7008.6113
$ws.Cells.Item(86, 9).Value = 5650.778
$ws.Cells.Item(86, 10).Value = 8366.444
$ws.Cells.Item(86, 11).Value = 5650.778
$ws.Cells.Item(86, 12).Value = 8366.444
$ws.Cells.Item(86, 13).Value = -4527.778
$ws.Cells.Item(86, 14).Value = -10612.444
$ws.Cells.Item(89, 8).Value = 7008.6113
$ws.Cells.Item(89, 9).Value = 5650.778
$ws.Cells.Item(89, 10).Value = 8366.444
$ws.Cells.Item(89, 11).Value = 28253.89
$ws.Cells.Item(89, 12).Value = 41832.22
$ws.Cells.Item(89, 13).Value = -22637.89
$ws.Cells.Item(89, 14).Value = -53064.22
$ws.Cells.Item(99, 8).Value = 5000
$ws.Cells.Item(99, 9).Value = 5000
$ws.Cells.Item(99, 11).Value = 5000
$ws.Cells.Item(99, 13).Value = -3502
$ws.Cells.Item(105, 8).Value = 2491.1
$ws.Cells.Item(105, 9).Value = 2442.0625
$ws.Cells.Item(105, 11).Value = 2442.0625
$ws.Cells.Item(105, 13).Value = -695.0625
$ws.Cells.Item(120, 8).Value = 45181.816
$ws.Cells.Item(120, 10).Value = 45200
$ws.Cells.Item(120, 12).Value = 45200
$ws.Cells.Item(120, 14).Value = -52458
$ws.Cells.Item(122, 8).Value = 1027.2307
$ws.Cells.Item(122, 9).Value = 1145
$ws.Cells.Item(122, 10).Value = 838.8
$ws.Cells.Item(122, 11).Value = 3435
$ws.Cells.Item(122, 12).Value = 2516.4
$ws.Cells.Item(122, 13).Value = -985
$ws.Cells.Item(122, 14).Value = -7416.4
$ws.Cells.Item(126, 8).Value = 5000
$ws.Cells.Item(126, 9).Value = 5000
$ws.Cells.Item(126, 11).Value = 15000
$ws.Cells.Item(126, 13).Value = -12530
$ws.Cells.Item(132, 8).Value = 41870.4
$ws.Cells.Item(132, 9).Value = 4510.1665
$ws.Cells.Item(132, 10).Value = 191311.33
$ws.Cells.Item(132, 11).Value = 13530.4995
$ws.Cells.Item(132, 12).Value = 573933.99
$ws.Cells.Item(132, 13).Value = -11000.4995
$ws.Cells.Item(132, 14).Value = -578993.99
$ws.Cells.Item(134, 8).Value = 5021.676
$ws.Cells.Item(134, 9).Value = 4136.0938
$ws.Cells.Item(134, 11).Value = 12408.2814
$ws.Cells.Item(134, 13).Value = -9873.2814
$ws.Cells.Item(136, 8).Value = 4549
$ws.Cells.Item(136, 9).Value = 3375.5454
$ws.Cells.Item(136, 10).Value = 6162.5
$ws.Cells.Item(136, 11).Value = 10126.6362
$ws.Cells.Item(136, 12).Value = 18487.5
$ws.Cells.Item(136, 13).Value = -7576.636200000001
$ws.Cells.Item(136, 14).Value = -23587.5
$ws.Cells.Item(141, 8).Value = 251756.27
$ws.Cells.Item(141, 9).Value = 100000
$ws.Cells.Item(141, 10).Value = 258982.77
$ws.Cells.Item(141, 11).Value = 100000
$ws.Cells.Item(141, 12).Value = 258982.77
$ws.Cells.Item(141, 13).Value = -94820
$ws.Cells.Item(141, 14).Value = -269342.77

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 3693.2632
$ws.Cells.Item(2, 10).Value = 6970.5
$ws.Cells.Item(2, 12).Value = 41823
$ws.Cells.Item(2, 14).Value = -42049
$ws.Cells.Item(4, 8).Value = 2349502.5
$ws.Cells.Item(4, 9).Value = 747983.9399999999
$ws.Cells.Item(4, 10).Value = 4980569
$ws.Cells.Item(4, 11).Value = 2243951.82
$ws.Cells.Item(4, 12).Value = 14941707
$ws.Cells.Item(4, 13).Value = -2243839.82
$ws.Cells.Item(4, 14).Value = -14941931
$ws.Cells.Item(5, 8).Value = 604.8
$ws.Cells.Item(5, 9).Value = 520.1111
$ws.Cells.Item(5, 10).Value = 674.0909
$ws.Cells.Item(5, 11).Value = 1560.3333
$ws.Cells.Item(5, 12).Value = 2022.2727
$ws.Cells.Item(5, 13).Value = -1448.3333
$ws.Cells.Item(5, 14).Value = -2246.2727
$ws.Cells.Item(6, 8).Value = 262.125
$ws.Cells.Item(6, 9).Value = 179.83333
$ws.Cells.Item(6, 11).Value = 539.49999
$ws.Cells.Item(6, 13).Value = -426.49999
$ws.Cells.Item(8, 8).Value = 316.36365
$ws.Cells.Item(8, 9).Value = 316.36365
$ws.Cells.Item(8, 11).Value = 949.09095
$ws.Cells.Item(8, 13).Value = -810.09095
$ws.Cells.Item(10, 8).Value = 386.22223
$ws.Cells.Item(10, 9).Value = 56
$ws.Cells.Item(10, 10).Value = 799
$ws.Cells.Item(10, 11).Value = 168
$ws.Cells.Item(10, 12).Value = 2397
$ws.Cells.Item(10, 13).Value = -29
$ws.Cells.Item(10, 14).Value = -2675
$ws.Cells.Item(12, 8).Value = 255.42857
$ws.Cells.Item(12, 10).Value = 255.42857
$ws.Cells.Item(12, 12).Value = 766.28571
$ws.Cells.Item(12, 14).Value = -1112.28571
$ws.Cells.Item(23, 8).Value = 960.4
$ws.Cells.Item(23, 9).Value = 300.5
$ws.Cells.Item(23, 11).Value = 901.5
$ws.Cells.Item(23, 13).Value = -666.5
$ws.Cells.Item(57, 8).Value = 2295.8333
$ws.Cells.Item(57, 10).Value = 2643.75
$ws.Cells.Item(57, 12).Value = 7931.25
$ws.Cells.Item(57, 14).Value = -9049.25
$ws.Cells.Item(68, 8).Value = 3830
$ws.Cells.Item(68, 10).Value = 4500
$ws.Cells.Item(68, 12).Value = 13500
$ws.Cells.Item(68, 14).Value = -15122
$ws.Cells.Item(71, 8).Value = 3830
$ws.Cells.Item(71, 10).Value = 4500
$ws.Cells.Item(71, 12).Value = 40500
$ws.Cells.Item(71, 14).Value = -48612
$ws.Cells.Item(98, 8).Value = 576.5454999999999
$ws.Cells.Item(98, 10).Value = 607.2
$ws.Cells.Item(98, 12).Value = 1821.6
$ws.Cells.Item(98, 14).Value = -4817.6
$ws.Cells.Item(107, 8).Value = 1009.5
$ws.Cells.Item(107, 9).Value = 960.6667
$ws.Cells.Item(107, 11).Value = 2882.0001
$ws.Cells.Item(107, 13).Value = -962.0001000000002
$ws.Cells.Item(112, 8).Value = 50005500
$ws.Cells.Item(112, 9).Value = 50005500
$ws.Cells.Item(112, 11).Value = 150016500
$ws.Cells.Item(112, 13).Value = -150015392
$ws.Cells.Item(135, 8).Value = 604.8
$ws.Cells.Item(135, 9).Value = 520.1111
$ws.Cells.Item(135, 10).Value = 674.0909
$ws.Cells.Item(135, 11).Value = 4680.9999
$ws.Cells.Item(135, 12).Value = 6066.8181
$ws.Cells.Item(135, 13).Value = -2145.9999
$ws.Cells.Item(135, 14).Value = -11136.8181
$ws.Cells.Item(137, 8).Value = 7603.6113
$ws.Cells.Item(137, 9).Value = 12502.444
$ws.Cells.Item(137, 10).Value = 2704.7778
$ws.Cells.Item(137, 11).Value = 37507.33199999999
$ws.Cells.Item(137, 12).Value = 8114.3334
$ws.Cells.Item(137, 13).Value = -32407.33199999999
$ws.Cells.Item(137, 14).Value = -18314.3334
$ws.Cells.Item(140, 8).Value = 19231706
$ws.Cells.Item(140, 9).Value = 20834012
$ws.Cells.Item(140, 10).Value = 4033
$ws.Cells.Item(140, 11).Value = 62502036
$ws.Cells.Item(140, 12).Value = 12099
$ws.Cells.Item(140, 13).Value = -62496856
$ws.Cells.Item(140, 14).Value = -22459
$ws.Cells.Item(141, 8).Value = 6815.6
$ws.Cells.Item(141, 9).Value = 6761.25
$ws.Cells.Item(141, 11).Value = 20283.75
$ws.Cells.Item(141, 13).Value = -15103.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 8365877
$ws.Cells.Item(11, 9).Value = 13333335
$ws.Cells.Item(11, 10).Value = 7011115.5
$ws.Cells.Item(11, 11).Value = 13333335
$ws.Cells.Item(11, 12).Value = 7011115.5
$ws.Cells.Item(11, 13).Value = -13333196
$ws.Cells.Item(11, 14).Value = -7011393.5
$ws.Cells.Item(20, 8).Value = 59090.453
$ws.Cells.Item(20, 9).Value = 200000
$ws.Cells.Item(20, 10).Value = 44999.5
$ws.Cells.Item(20, 11).Value = 200000
$ws.Cells.Item(20, 12).Value = 44999.5
$ws.Cells.Item(20, 13).Value = -199755
$ws.Cells.Item(20, 14).Value = -45489.5
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).ClearContents()
$ws.Cells.Item(70, 14).ClearContents()
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).ClearContents()
$ws.Cells.Item(73, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 3036.182
$ws.Cells.Item(80, 9).Value = 2339.9
$ws.Cells.Item(80, 11).Value = 2339.9
$ws.Cells.Item(80, 13).Value = -1341.9
$ws.Cells.Item(83, 8).Value = 3036.182
$ws.Cells.Item(83, 9).Value = 2339.9
$ws.Cells.Item(83, 11).Value = 11699.5
$ws.Cells.Item(83, 13).Value = -6707.5
$ws.Cells.Item(93, 8).Value = 39733
$ws.Cells.Item(93, 10).Value = 39733
$ws.Cells.Item(93, 12).Value = 39733
$ws.Cells.Item(93, 14).Value = -43477
$ws.Cells.Item(97, 8).Value = 1382.4857
$ws.Cells.Item(97, 9).Value = 1480.12
$ws.Cells.Item(97, 11).Value = 1480.12
$ws.Cells.Item(97, 13).Value = -984.1199999999999
$ws.Cells.Item(113, 8).Value = 7029.7144
$ws.Cells.Item(113, 9).Value = 4423.75
$ws.Cells.Item(113, 10).Value = 10504.333
$ws.Cells.Item(113, 11).Value = 4423.75
$ws.Cells.Item(113, 12).Value = 10504.333
$ws.Cells.Item(113, 13).Value = -2253.75
$ws.Cells.Item(113, 14).Value = -14844.333
$ws.Cells.Item(114, 8).Value = 65927.5
$ws.Cells.Item(114, 10).Value = 79855
$ws.Cells.Item(114, 12).Value = 79855
$ws.Cells.Item(114, 14).Value = -88533
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 2091.6365
$ws.Cells.Item(122, 9).Value = 1888.8889
$ws.Cells.Item(122, 10).Value = 3004
$ws.Cells.Item(122, 11).Value = 5666.6667
$ws.Cells.Item(122, 12).Value = 9012
$ws.Cells.Item(122, 13).Value = -3216.6667
$ws.Cells.Item(122, 14).Value = -13912
$ws.Cells.Item(126, 8).Value = 3992.9167
$ws.Cells.Item(126, 9).Value = 3122.6667
$ws.Cells.Item(126, 10).Value = 4863.1665
$ws.Cells.Item(126, 11).Value = 9368.000100000001
$ws.Cells.Item(126, 12).Value = 14589.4995
$ws.Cells.Item(126, 13).Value = -6898.000100000001
$ws.Cells.Item(126, 14).Value = -19529.4995
$ws.Cells.Item(132, 8).Value = 4500.436
$ws.Cells.Item(132, 9).Value = 3374.0293
$ws.Cells.Item(132, 10).Value = 12160
$ws.Cells.Item(132, 11).Value = 10122.0879
$ws.Cells.Item(132, 12).Value = 36480
$ws.Cells.Item(132, 13).Value = -7592.0879
$ws.Cells.Item(132, 14).Value = -41540
$ws.Cells.Item(136, 8).Value = 98714.625
$ws.Cells.Item(136, 10).Value = 98714.625
$ws.Cells.Item(136, 12).Value = 296143.875
$ws.Cells.Item(136, 14).Value = -301243.875
$ws.Cells.Item(140, 8).Value = 119999.5
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 119999.5
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 119999.5
$ws.Cells.Item(140, 13).ClearContents()
$ws.Cells.Item(140, 14).Value = -130359.5
$ws.Cells.Item(141, 8).Value = 96631.664
$ws.Cells.Item(141, 10).Value = 96631.664
$ws.Cells.Item(141, 12).Value = 96631.664
$ws.Cells.Item(141, 14).Value = -106991.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 10851.5
$ws.Cells.Item(3, 9).Value = 100
$ws.Cells.Item(3, 10).Value = 13001.8
$ws.Cells.Item(3, 11).Value = 100
$ws.Cells.Item(3, 12).Value = 13001.8
$ws.Cells.Item(3, 13).Value = 12
$ws.Cells.Item(3, 14).Value = -13225.8
$ws.Cells.Item(7, 8).Value = 2831.5
$ws.Cells.Item(7, 9).Value = 2831.5
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 2831.5
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -2719.5
$ws.Cells.Item(7, 14).ClearContents()
$ws.Cells.Item(10, 8).Value = 4634.6665
$ws.Cells.Item(10, 9).Value = 850
$ws.Cells.Item(10, 10).Value = 6527
$ws.Cells.Item(10, 11).Value = 850
$ws.Cells.Item(10, 12).Value = 6527
$ws.Cells.Item(10, 13).Value = -710
$ws.Cells.Item(10, 14).Value = -6807
$ws.Cells.Item(15, 8).Value = 10851.5
$ws.Cells.Item(15, 9).Value = 100
$ws.Cells.Item(15, 10).Value = 13001.8
$ws.Cells.Item(15, 11).Value = 100
$ws.Cells.Item(15, 12).Value = 13001.8
$ws.Cells.Item(15, 13).Value = 70
$ws.Cells.Item(15, 14).Value = -13341.8
$ws.Cells.Item(16, 8).Value = 1560.0555
$ws.Cells.Item(16, 9).Value = 1269.4706
$ws.Cells.Item(16, 10).Value = 6500
$ws.Cells.Item(16, 11).Value = 1269.4706
$ws.Cells.Item(16, 12).Value = 6500
$ws.Cells.Item(16, 13).Value = -1099.4706
$ws.Cells.Item(16, 14).Value = -6840
$ws.Cells.Item(17, 8).Value = 15000
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 14).ClearContents()
$ws.Cells.Item(25, 8).Value = 12500
$ws.Cells.Item(25, 9).Value = 15000
$ws.Cells.Item(25, 11).Value = 15000
$ws.Cells.Item(25, 13).Value = -14770
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 14).ClearContents()
$ws.Cells.Item(46, 8).Value = 5586.091
$ws.Cells.Item(46, 9).Value = 1321.7778
$ws.Cells.Item(46, 10).Value = 8538.308000000001
$ws.Cells.Item(46, 11).Value = 1321.7778
$ws.Cells.Item(46, 12).Value = 8538.308000000001
$ws.Cells.Item(46, 13).Value = -1133.7778
$ws.Cells.Item(46, 14).Value = -8914.308000000001
$ws.Cells.Item(61, 8).Value = 7991.8335
$ws.Cells.Item(61, 9).Value = 7991.8335
$ws.Cells.Item(61, 11).Value = 7991.8335
$ws.Cells.Item(61, 13).Value = -7789.8335
$ws.Cells.Item(68, 8).Value = 3749.75
$ws.Cells.Item(68, 9).Value = 2571.1428
$ws.Cells.Item(68, 11).Value = 2571.1428
$ws.Cells.Item(68, 13).Value = -1822.1428
$ws.Cells.Item(70, 8).Value = 80000
$ws.Cells.Item(70, 10).Value = 80000
$ws.Cells.Item(70, 12).Value = 80000
$ws.Cells.Item(70, 14).Value = -80540
$ws.Cells.Item(71, 8).Value = 3749.75
$ws.Cells.Item(71, 9).Value = 2571.1428
$ws.Cells.Item(71, 11).Value = 12855.714
$ws.Cells.Item(71, 13).Value = -9111.714
$ws.Cells.Item(73, 8).Value = 80000
$ws.Cells.Item(73, 10).Value = 80000
$ws.Cells.Item(73, 12).Value = 80000
$ws.Cells.Item(73, 14).Value = -81872
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 14).ClearContents()
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 14).ClearContents()
$ws.Cells.Item(93, 8).Value = 3045.7
$ws.Cells.Item(93, 9).Value = 2842.2942
$ws.Cells.Item(93, 10).Value = 4198.3335
$ws.Cells.Item(93, 11).Value = 2842.2942
$ws.Cells.Item(93, 12).Value = 4198.3335
$ws.Cells.Item(93, 13).Value = -1594.2942
$ws.Cells.Item(93, 14).Value = -6694.3335
$ws.Cells.Item(100, 8).Value = 5557543.5
$ws.Cells.Item(100, 9).Value = 8334836.5
$ws.Cells.Item(100, 10).Value = 2958.2
$ws.Cells.Item(100, 11).Value = 8334836.5
$ws.Cells.Item(100, 12).Value = 2958.2
$ws.Cells.Item(100, 13).Value = -8334295.5
$ws.Cells.Item(100, 14).Value = -4040.2
$ws.Cells.Item(113, 8).Value = 7991.8335
$ws.Cells.Item(113, 9).Value = 7991.8335
$ws.Cells.Item(113, 11).Value = 7991.8335
$ws.Cells.Item(113, 13).Value = -5821.8335
$ws.Cells.Item(122, 8).Value = 29478612
$ws.Cells.Item(122, 9).Value = 41674540
$ws.Cells.Item(122, 11).Value = 125023620
$ws.Cells.Item(122, 13).Value = -125021170
$ws.Cells.Item(126, 8).Value = 2831.5
$ws.Cells.Item(126, 9).Value = 2831.5
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 8494.5
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -6024.5
$ws.Cells.Item(126, 14).ClearContents()
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 7871.3267
$ws.Cells.Item(132, 9).Value = 7210.7812
$ws.Cells.Item(132, 10).Value = 9114.706
$ws.Cells.Item(132, 11).Value = 21632.3436
$ws.Cells.Item(132, 12).Value = 27344.118
$ws.Cells.Item(132, 13).Value = -19102.3436
$ws.Cells.Item(132, 14).Value = -32404.118
$ws.Cells.Item(136, 8).Value = 6995.9585
$ws.Cells.Item(136, 9).Value = 5350.1577
$ws.Cells.Item(136, 10).Value = 13250
$ws.Cells.Item(136, 11).Value = 16050.4731
$ws.Cells.Item(136, 12).Value = 39750
$ws.Cells.Item(136, 13).Value = -13500.4731
$ws.Cells.Item(136, 14).Value = -44850
$ws.Cells.Item(140, 8).Value = 536998.6
$ws.Cells.Item(140, 10).Value = 536998.6
$ws.Cells.Item(140, 12).Value = 536998.6
$ws.Cells.Item(140, 14).Value = -547358.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(20, 8).Value = 13500
$ws.Cells.Item(20, 9).Value = 13000
$ws.Cells.Item(20, 11).Value = 13000
$ws.Cells.Item(20, 13).Value = -12760
$ws.Cells.Item(42, 8).Value = 40000
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 13).ClearContents()
$ws.Cells.Item(68, 8).Value = 20070.5
$ws.Cells.Item(68, 10).Value = 20070.5
$ws.Cells.Item(68, 12).Value = 20070.5
$ws.Cells.Item(68, 14).Value = -21692.5
$ws.Cells.Item(71, 8).Value = 20070.5
$ws.Cells.Item(71, 10).Value = 20070.5
$ws.Cells.Item(71, 12).Value = 60211.5
$ws.Cells.Item(71, 14).Value = -68323.5
$ws.Cells.Item(92, 8).Value = 31499.6
$ws.Cells.Item(92, 10).Value = 31499.6
$ws.Cells.Item(92, 12).Value = 31499.6
$ws.Cells.Item(92, 14).Value = -36491.6
$ws.Cells.Item(107, 8).Value = 674.4706
$ws.Cells.Item(107, 9).Value = 656.4
$ws.Cells.Item(107, 11).Value = 1969.2
$ws.Cells.Item(107, 13).Value = -49.19999999999982
$ws.Cells.Item(117, 8).Value = 50000
$ws.Cells.Item(117, 10).Value = 50000
$ws.Cells.Item(117, 12).Value = 50000
$ws.Cells.Item(117, 14).Value = -59178
$ws.Cells.Item(122, 8).Value = 5932.0713
$ws.Cells.Item(122, 9).Value = 5914.636
$ws.Cells.Item(122, 10).Value = 5996
$ws.Cells.Item(122, 11).Value = 17743.908
$ws.Cells.Item(122, 12).Value = 17988
$ws.Cells.Item(122, 13).Value = -15293.908
$ws.Cells.Item(122, 14).Value = -22888
$ws.Cells.Item(132, 8).Value = 4238.5527
$ws.Cells.Item(132, 9).Value = 3364.5173
$ws.Cells.Item(132, 11).Value = 10093.5519
$ws.Cells.Item(132, 13).Value = -7563.5519
$ws.Cells.Item(136, 8).Value = 2277.5894
$ws.Cells.Item(136, 9).Value = 1488.8718
$ws.Cells.Item(136, 11).Value = 4466.6154
$ws.Cells.Item(136, 13).Value = -1916.6154
$ws.Cells.Item(137, 8).Value = 119999
$ws.Cells.Item(137, 10).Value = 119999
$ws.Cells.Item(137, 12).Value = 119999
$ws.Cells.Item(137, 14).Value = -130199

